$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark left over from the previous
#    edit session (sits right after the Cloudfront sentence).
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) SWC cert date typo: "Dec 2019" -> "Dec 2018"
#    Only the final digit changes; split it into its own run (as
#    Word would when you retype a single character) by nudging the
#    character formatting and then restoring it.
# ------------------------------------------------------------------
$swc = $d.Content.Duplicate
$swc.Find.Execute("Dec 2019") | Out-Null
$swcEnd = $swc.End

$lastDigit = $d.Range($swcEnd - 1, $swcEnd)
$lastDigit.Font.Bold = 1
$lastDigit.Text = "8"
$lastDigit2 = $d.Range($swcEnd - 1, $swcEnd)
$lastDigit2.Font.Bold = 0

# ------------------------------------------------------------------
# 3) SumoLogic cert casing typo: "MAY 2019" -> "May 2019"
#    Split "AY" into its own run (becoming "ay") and drop the
#    "_GoBack" bookmark at the point right after it, matching where
#    Word would leave it after the last edit.
# ------------------------------------------------------------------
$sumo = $d.Content.Duplicate
$sumo.Find.Execute("Pro User and Power User Certified, MAY 2019") | Out-Null
$sumoStart = $sumo.Start

$ayStart = $sumoStart + 36
$ayEnd = $sumoStart + 38

$ay = $d.Range($ayStart, $ayEnd)
$ay.Font.Bold = 1
$ay.Text = "ay"
$ay2 = $d.Range($ayStart, $ayStart + 2)
$ay2.Font.Bold = 0

$goBackPoint = $d.Range($ayStart + 2, $ayStart + 2)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
